$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates from the refreshed crypto price feed.
# D-column "Price" cells are forced to Text so numeric-looking values
# (e.g. "4.13", "0.0950") are not auto-converted/rounded by Excel,
# matching the inlineStr storage used by the source workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.408.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.570.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.16"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -0.79%  "
$ws.Range("E10").Value = "  +0.56%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.794.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.565.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.386.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0691"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "212.72"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0471"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.392.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.943"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.528"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E42").Value = "  +2.67%  "
$ws.Range("E43").Value = "  +4.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.706.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.88%  "
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0950"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.93%  "
